$wb = $excel.ActiveWorkbook

# Sheet ALC, row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 181.4762
$ws.Range("I15").Value = 181.4762
$ws.Range("K15").Value = 544.4286
$ws.Range("M15").Value = -375.4286

# Sheet ALC, row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1716.1666
$ws.Range("I38").Value = 65.59999999999999
$ws.Range("K38").Value = 196.8
$ws.Range("M38").Value = 175.2

# Sheet ALC, row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6694.2
$ws.Range("J51").Value = 6694.2
$ws.Range("L51").Value = 6694.2
$ws.Range("N51").Value = -7662.2

# Sheet ALC, row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3430.5
$ws.Range("I98").Value = 362
$ws.Range("J98").Value = 6499
$ws.Range("K98").Value = 362
$ws.Range("L98").Value = 6499
$ws.Range("M98").Value = 1136
$ws.Range("N98").Value = -9495

# Sheet ALC, row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3430.5
$ws.Range("I122").Value = 362
$ws.Range("J122").Value = 6499
$ws.Range("K122").Value = 1086
$ws.Range("L122").Value = 19497
$ws.Range("M122").Value = 1364
$ws.Range("N122").Value = -24397

# Sheet ALC, row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 5399.6665
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 5399.6665
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 48596.9985
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -53516.9985

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5279.75
$ws.Range("I132").Value = 5279.75
$ws.Range("K132").Value = 15839.25
$ws.Range("M132").Value = -13309.25

# Sheet ALC, row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 662.25
$ws.Range("I135").Value = 662.25
$ws.Range("K135").Value = 5960.25
$ws.Range("M135").Value = -3425.25

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3080.8572
$ws.Range("I2").Value = 3080.8572
$ws.Range("K2").Value = 3080.8572
$ws.Range("M2").Value = -2967.8572

# Sheet ARM, row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 251.375
$ws.Range("I4").Value = 215.85715
$ws.Range("K4").Value = 215.85715
$ws.Range("M4").Value = -99.85714999999999

# Sheet ARM, row 14
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 4250
$ws.Range("I14").Value = 4500
$ws.Range("J14").Value = 4000
$ws.Range("K14").Value = 4500
$ws.Range("L14").Value = 4000
$ws.Range("M14").Value = -4325
$ws.Range("N14").Value = -4350

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 511
$ws.Range("I32").Value = 511
$ws.Range("K32").Value = 511
$ws.Range("M32").Value = -224

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1531.4286
$ws.Range("I61").Value = 1531.4286
$ws.Range("K61").Value = 1531.4286
$ws.Range("M61").Value = -1319.4286

# Sheet ARM, row 95
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3080.8572
$ws.Range("I116").Value = 3080.8572
$ws.Range("K116").Value = 3080.8572
$ws.Range("M116").Value = -786.8571999999999

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1531.4286
$ws.Range("I136").Value = 1531.4286
$ws.Range("K136").Value = 4594.2858
$ws.Range("M136").Value = -2044.2858

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3080.8572
$ws.Range("I3").Value = 3080.8572
$ws.Range("K3").Value = 3080.8572
$ws.Range("M3").Value = -2966.8572

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2357.75
$ws.Range("I134").Value = 2357.75
$ws.Range("K134").Value = 7073.25
$ws.Range("M134").Value = -4538.25

# Sheet CRP, row 69
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 20666.666
$ws.Range("I69").Value = 2000
$ws.Range("K69").Value = 2000
$ws.Range("M69").Value = -1251

# Sheet CRP, row 72
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 20666.666
$ws.Range("I72").Value = 2000
$ws.Range("K72").Value = 6000
$ws.Range("M72").Value = -2256

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1654.1538
$ws.Range("I132").Value = 1375.3334
$ws.Range("K132").Value = 4126.0002
$ws.Range("M132").Value = -1596.0002

# Sheet CUL, row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 77.933334
$ws.Range("I2").Value = 55.8
$ws.Range("J2").Value = 122.2
$ws.Range("K2").Value = 334.8
$ws.Range("L2").Value = 733.2
$ws.Range("M2").Value = -221.8
$ws.Range("N2").Value = -959.2

# Sheet CUL, row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5168166.5
$ws.Range("I4").Value = 336333.34
$ws.Range("J4").Value = 10000000
$ws.Range("K4").Value = 1009000.02
$ws.Range("L4").Value = 30000000
$ws.Range("M4").Value = -1008888.02
$ws.Range("N4").Value = -30000224

# Sheet CUL, row 8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 913.6667
$ws.Range("I8").Value = 913.6667
$ws.Range("K8").Value = 2741.0001
$ws.Range("M8").Value = -2602.0001

# Sheet CUL, row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2000
$ws.Range("J55").Value = 2000
$ws.Range("L55").Value = 6000
$ws.Range("N55").Value = -6354

# Sheet GSM, row 6
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 1550
$ws.Range("I6").Value = 600
$ws.Range("J6").Value = 2500
$ws.Range("K6").Value = 600
$ws.Range("L6").Value = 2500
$ws.Range("M6").Value = -487
$ws.Range("N6").Value = -2726

# Sheet GSM, row 16
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H16").Value = 1550
$ws.Range("I16").Value = 600
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 600
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -350

# Sheet GSM, row 44
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()

# Sheet GSM, row 58
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 40000
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

# Sheet GSM, row 103
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 249999.5
$ws.Range("J103").Value = 249999.5
$ws.Range("L103").Value = 249999.5
$ws.Range("N103").Value = -252343.5

# Sheet GSM, row 136
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 25000
$ws.Range("J136").Value = 25000
$ws.Range("L136").Value = 75000
$ws.Range("N136").Value = -80100

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2889.8
$ws.Range("I46").Value = 2789
$ws.Range("K46").Value = 2789
$ws.Range("M46").Value = -2601

# Sheet WVR, row 7
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 3500
$ws.Range("I7").Value = 2750
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 2750
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -2637
$ws.Range("N7").Value = -5226

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3525.8
$ws.Range("I132").Value = 3084.2222
$ws.Range("K132").Value = 9252.6666
$ws.Range("M132").Value = -6722.6666

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4912
$ws.Range("I136").Value = 4912
$ws.Range("K136").Value = 14736
$ws.Range("M136").Value = -12186
